$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.845.17"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "3.512.46"
$ws.Range("E3").Value = "  -3.61%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "193.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.05%  "
$ws.Range("D7").Value = "3.500.66"
$ws.Range("E7").Value = "  -3.68%  "
$ws.Range("E8").Value = "  -2.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  -6.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.623"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.16%  "
$ws.Range("E12").Value = "  -4.27%  "
$ws.Range("E13").Value = "  -6.58%  "
$ws.Range("E14").Value = "  -4.43%  "
$ws.Range("D15").Value = "4.075.69"
$ws.Range("E15").Value = "  -3.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "645.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.89%  "
$ws.Range("D17").Value = "69.825.44"
$ws.Range("E17").Value = "  -1.80%  "
$ws.Range("D18").Value = "3.513.12"
$ws.Range("E18").Value = "  -3.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.50%  "
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.953"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.99%  "
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.29%  "
$ws.Range("E26").Value = "  -7.14%  "
$ws.Range("E27").Value = "  -4.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.68%  "
$ws.Range("E33").Value = "  -4.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "577.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.02%  "
$ws.Range("E35").Value = "  -4.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "61.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.78%  "
$ws.Range("D37").Value = "3.781.60"
$ws.Range("E37").Value = "  -3.66%  "
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0797"
$ws.Range("E39").Value = "  -8.39%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("B41").Value = "CoreDAO"
$ws.Range("C41").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +39.53%  "
$ws.Range("E42").Value = "  -4.21%  "
$ws.Range("E43").Value = "  -4.08%  "
$ws.Range("E44").Value = "  -2.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "34.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.18%  "
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("E47").Value = "  -5.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.20%  "
$ws.Range("E49").Value = "  -3.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("E51").Value = "  -5.09%  "
